$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the duplicate number-format style (style index 12) so these
# cells share the same style as the rest of the "00" formatted cells
# (style index 1). Doing this cell-by-cell (rather than as one big range)
# avoids touching the untouched F/H columns in the same rows.
$legacyStyledCells = @(
    "E23","G23",
    "E24","G24","I24",
    "E25","G25","I25",
    "E26","G26",
    "E27","G27",
    "E28","G28",
    "E29","G29","I29",
    "E30","G30",
    "G31","I31",
    "G32","I32",
    "G33","I33",
    "G34","I34"
)
foreach ($addr in $legacyStyledCells) {
    $ws.Range($addr).NumberFormat = "00"
}

# --- Complete Ch 08 05 "Understanding semi-additive calculations":
# append the two new rows describing LASTDATE/LASTNONBLANK and opening
# and closing balances.
$ws.Range("E37").Value = 5
$ws.Range("F37").Value = "Understanding semi-additive calculations"
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = "Using LASTDATE and LASTNONBLANK"
$ws.Range("I37").Value = "LASTDATE, LASTNONBLANK"

$ws.Range("E38").Value = 5
$ws.Range("F38").Value = "Understanding semi-additive calculations"
$ws.Range("G38").Value = 2
$ws.Range("H38").Value = "Working with opening and closing balances"

# Match the "00" number format used elsewhere in the table for the new
# numeric section/subsection index cells.
$newNumberCells = @("E37","G37","I37","E38","G38")
foreach ($addr in $newNumberCells) {
    $ws.Range($addr).NumberFormat = "00"
}

# --- Update the view/selection to match where the author left off editing.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F39").Select()
